# Update the "date" column (C) values on the active sheet.
# Each entry maps a row number to its new date string. The leading
# apostrophe forces Excel to store the value as literal text instead of
# auto-converting the date-like string into a date serial number, which
# matches the original inline-string cell content (e.g. "2023-11-01").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDates = @(
    @{ Row = 5;  Date = "2023-11-01" },
    @{ Row = 6;  Date = "2023-11-01" },
    @{ Row = 7;  Date = "2023-11-01" },
    @{ Row = 8;  Date = "2023-11-02" },
    @{ Row = 9;  Date = "2023-11-02" },
    @{ Row = 10; Date = "2023-11-02" },
    @{ Row = 11; Date = "2023-11-02" },
    @{ Row = 12; Date = "2023-11-03" },
    @{ Row = 13; Date = "2023-11-03" },
    @{ Row = 14; Date = "2023-11-04" },
    @{ Row = 15; Date = "2023-11-05" },
    @{ Row = 16; Date = "2023-11-06" },
    @{ Row = 17; Date = "2023-11-06" },
    @{ Row = 18; Date = "2023-11-07" },
    @{ Row = 19; Date = "2023-11-07" },
    @{ Row = 20; Date = "2023-11-08" },
    @{ Row = 21; Date = "2023-11-08" },
    @{ Row = 22; Date = "2023-11-09" },
    @{ Row = 23; Date = "2023-11-09" }
)

foreach ($entry in $newDates) {
    $ws.Cells.Item($entry.Row, 3).Value = "'" + $entry.Date
}
